# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.434.87"
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").Value = "3.006.93"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'508.24"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").Value = "'138.94"
$ws.Range("E6").Value = "  +1.05%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.437"
$ws.Range("E8").Value = "  +0.89%  "
$ws.Range("D9").Value = "'7.50"
$ws.Range("E9").Value = "  -1.68%  "
$ws.Range("E10").Value = "  +1.43%  "
$ws.Range("D11").Value = "'0.365"
$ws.Range("E11").Value = "  +3.30%  "
$ws.Range("D12").Value = "3.521.47"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("D14").Value = "'26.33"
$ws.Range("E14").Value = "  +3.17%  "
$ws.Range("D15").Value = "'0.0000163"
$ws.Range("E15").Value = "  +5.54%  "
$ws.Range("D16").Value = "57.461.92"
$ws.Range("E16").Value = "  +1.46%  "
$ws.Range("D17").Value = "'6.19"
$ws.Range("E17").Value = "  +6.25%  "
$ws.Range("D18").Value = "3.010.21"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").Value = "'12.79"
$ws.Range("E19").Value = "  +1.93%  "
$ws.Range("D20").Value = "'7.95"
$ws.Range("E20").Value = "  +1.14%  "
$ws.Range("D21").Value = "'329.81"
$ws.Range("E21").Value = "  +0.44%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "'0.497"
$ws.Range("E23").Value = "  +3.55%  "
$ws.Range("D24").Value = "'64.45"
$ws.Range("E24").Value = "  +3.08%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "0.0₃0916"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").Value = "'6.79"
$ws.Range("E28").Value = "  +3.00%  "
$ws.Range("D29").Value = "'7.39"
$ws.Range("E29").Value = "  +5.56%  "
$ws.Range("E30").Value = "  +2.13%  "
$ws.Range("D31").Value = "'1.18"
$ws.Range("E31").Value = "  -5.61%  "
$ws.Range("D32").Value = "'20.56"
$ws.Range("E32").Value = "  -0.91%  "
$ws.Range("D33").Value = "'4.71"
$ws.Range("E33").Value = "  +4.47%  "
$ws.Range("D34").Value = "'154.08"
$ws.Range("E34").Value = "  -1.22%  "
$ws.Range("D35").Value = "'5.87"
$ws.Range("E35").Value = "  +4.10%  "
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.0682"
$ws.Range("E37").Value = "  +1.23%  "
$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D38").Value = "'24.37"
$ws.Range("E38").Value = "  +0.32%  "
$ws.Range("D39").Value = "3.040.01"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("D40").Value = "'37.35"
$ws.Range("E40").Value = "  +1.88%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'3.84"
$ws.Range("E42").Value = "  +5.97%  "
$ws.Range("D43").Value = "2.278.22"
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("D44").Value = "'0.648"
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("E45").Value = "  -0.26%  "
$ws.Range("D46").Value = "'0.983"
$ws.Range("E46").Value = "  -2.04%  "
$ws.Range("D47").Value = "'6.00"
$ws.Range("E47").Value = "  +3.78%  "
$ws.Range("E48").Value = "  +1.31%  "
$ws.Range("D49").Value = "'19.37"
$ws.Range("E49").Value = "  +1.08%  "
$ws.Range("D50").Value = "'1.85"
$ws.Range("E50").Value = "  -7.26%  "
$ws.Range("D51").Value = "'0.0893"
$ws.Range("E51").Value = "  +2.06%  "
